# Changed time display in array to seconds
# The "TotalTime" column (H) on each sheet held values as text (milliseconds).
# They are converted to literal numeric values in seconds (divided by 1000).

$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("7-Level")
$ws5 = $wb.Worksheets.Item("5-Level")
$ws3 = $wb.Worksheets.Item("3-Level")

# --- 7-Level sheet ---
$ws7.Range("H3").Value2 = 0.0013272488609979
$ws7.Range("H4").Value2 = 0.0012733105671904001
$ws7.Range("H5").Value2 = 0.0013272488609979
$ws7.Range("H8").Value2 = 0.0012999999999999999
$ws7.Range("H10").Value2 = 0.0022603911202810998
$ws7.Range("H12").Value2 = 0.0021093738732015001
$ws7.Range("H14").Value2 = 0.0013272488609979
$ws7.Range("H16").Value2 = 0.0012733105671904001
$ws7.Range("H18").Value2 = 0.0013090202209187999
$ws7.Range("H20").Value2 = 0.0018117050564913999
$ws7.Range("H22").Value2 = 0.0012821455999707

# --- 5-Level sheet ---
$ws5.Range("H3").Value2 = 0.0013181030298264999
$ws5.Range("H4").Value2 = 0.0013090202209187999
$ws5.Range("H6").Value2 = 0.0012999999999999999
$ws5.Range("H7").Value2 = 0.0013090202209187999
$ws5.Range("H9").Value2 = 0.0013090202209187999
$ws5.Range("H11").Value2 = 0.0013090202209187999
$ws5.Range("H13").Value2 = 0.0013090202209187999

# --- 3-Level sheet ---
$ws3.Range("H3").Value2 = 0.0012821455999707
$ws3.Range("H4").Value2 = 0.0013181030298264999
$ws3.Range("H6").Value2 = 0.0012821455999707

# --- Selections / active sheet ---
# 5-Level: selection moves from B17:F17 (active F17) to H14, and the
# topLeftCell scroll anchor is cleared.
$ws5.Activate()
$ws5.Range("H14").Select()

# 7-Level: selection stays at K11 (tabSelected flag is removed from this sheet).
$ws7.Activate()
$ws7.Range("K11").Select()

# 3-Level becomes the selected/active tab, with selection moved to H7.
$ws3.Activate()
$ws3.Range("H7").Select()
